$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("countries")

# Fix the "Balken(e)de" typo to the correct "Balkenende" spelling in the
# cabinet-reference column (C) of the countries sheet.
$ws.Range("C147:C182").Value = "Balkenende I"
$ws.Range("C183:C218").Value = "Balkenende II"
$ws.Range("C219:C254").Value = "Balkenende III"
$ws.Range("C255:C290").Value = "Balkenende IV"

# Update the view state to reflect where the author was working: select the
# corrected "Balkenende IV" range, with C255 as the active cell.
$ws.Activate()
$ws.Range("C255:C290").Select()
